$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A13").Value = "Saturday 5.5.18"
$ws.Range("B13").Value = "1600-1700"
$ws.Range("C13").Value = 1

$ws.Range("D13").Select()
